$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate formatting/types of the last existing row (row 25) into the new row 26,
# so the date style and "TRUE"/"Maria" text cells keep the same representation
# (shared-string text rather than boolean/number) as the rest of the table.
$ws.Range("A25:E25").Copy($ws.Range("A26:E26"))

# Fill in the new submission's specific values.
# Set "Name Ramp" (C) before "Name" (B) so the new shared strings are appended
# to sharedStrings.xml in the same order as in the target workbook.
# (D26 "Hand in" and E26 "By" already carry the correct text values - TRUE / Maria -
# copied from row 25 above, so they are left untouched to avoid Excel re-coercing
# the literal string "TRUE" back into a boolean cell.)
$ws.Range("A26").Value = 44891
$ws.Range("C26").Value = "78LaundryIsDone"
$ws.Range("B26").Value = "221126_cat_v1data"

# Grow the table (ListObject) to include the new row.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:E26"))

# Update the active selection to reflect where the user would click next.
$ws.Range("B27").Select()
